# ActivationFluxConfiguration.xlsx - "agregar planes de residencial" edit
#
# On the "Plans" sheet, row 5 (plan #1, Residencial / Sin_TotalPlay_TV) is
# kept, its price (column D) is bumped from 50 to 100, and the three
# duplicate plan rows below it (rows 6-8, which repeated the same
# Residencial / Sin_TotalPlay_TV plan at different prices) are cleared out
# back to blank template rows, ready for new plan data to be entered.

$wb = $excel.ActiveWorkbook

# --- Data edits on the "Plans" sheet ---------------------------------
$plans = $wb.Worksheets.Item("Plans")

# Row 5: price 50 -> 100
$plans.Range("D5").Value = 100

# Rows 6-8: wipe out the Plan/SubPlan/Price columns (B:D) so the rows go
# back to being empty, unused template rows.
$plans.Range("B6:D8").ClearContents()

# --- Selection bookkeeping --------------------------------------------
# While making the edit above, the B6:D8 block ended up selected; that
# selection is also reflected on the "Introduction" sheet.
$intro = $wb.Worksheets.Item("Introduction")
$intro.Activate() | Out-Null
$intro.Range("B6:D8").Select() | Out-Null

# Leave "TestConfiguration" and "Tables" selections as they were.

# Finish back on "Plans" (the sheet that was actually edited) with the
# cleared B6:D8 block selected, matching the in-progress edit state.
$plans.Activate() | Out-Null
$plans.Range("B6:D8").Select() | Out-Null
